# Add an auto-selected "model" / complexity field to the Misc sheet,
# with a dropdown (data validation list) of the available model types.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Misc")

$ws.Range("A2").Value = "model"
$ws.Range("B2").Value = "Linear"

$validationRange = $ws.Range("B2")
$validationRange.Validation.Delete()
$validationRange.Validation.Add(3, 1, 1, '"Linear, 2FI, Quadratic"')
$validationRange.Validation.IgnoreBlank = $true
$validationRange.Validation.InCellDropdown = $true
$validationRange.Validation.ShowInput = $true
$validationRange.Validation.ShowError = $true

$ws.Range("H5").Select()
